# GeneracionSiniestroMotor.xlsx - "Completar actividades para realizar pagos.
# Verificar pagos en CC y en SISE"
#
# Row 4 (claim #3) is updated: the policy/claim reference number and the
# claim date are refreshed, and the severity/fast-track/description fields
# are completed to reflect a partial-loss claim, matching the pattern
# already used by the other rows in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NroPoliza (F4) and FechaSiniestro (H4) both look like numbers/dates, but
# must stay stored as plain text (as the rest of the sheet does), so a
# leading apostrophe is used to force text entry instead of numeric/date
# auto-conversion.
$ws.Range("H4").Value = "'07/04/2021"
$ws.Range("F4").Value = "'04104016054"

# Gravedad / EsFastTrack / Descripcion
$ws.Range("S4").Value = "Pérdidas parciales"
$ws.Range("T4").Value = "No"
$ws.Range("U4").Value = "Parcial (Resto del Vehiculo)"

# Reflect the author's final cursor position / scroll on the sheet.
$ws.Range("E3").Select()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
